$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=34.111822; H=102.335466; I=0.4228853893909983; J=0.4228853893909983; K=3; M=102.8289443333334; N=308.486833; O=0.5559120396302444; P=0.5559120396302443; Q=3507.682645546576; R=31569.14380991918; S=0.23508707934618; T=0.23508707934618 }
    3  = @{ E=3; G=34.111822; H=102.335466; I=0.4228853893909983; J=0.4228853893909983; K=3; M=63.66262833333334; N=190.987885; O=0.3441717873742006; P=0.3441717873742006; Q=2171.648245758823; R=19544.83421182941; S=0.1455452203211347; T=0.1455452203211347 }
    4  = @{ E=3; G=34.111822; H=102.335466; I=0.4228853893909983; J=0.4228853893909983; K=3; M=18.481835; N=55.445505; O=0.09991617299555507; P=0.09991617299555505; Q=630.44906575337; R=5674.041591780329; S=0.04225308972368366; T=0.04225308972368366 }
    5  = @{ E=3; G=34.88211266666666; H=104.646338; I=0.4324347083490296; J=0.4324347083490295; K=3; M=102.8289443333334; N=308.486833; O=0.5559120396302444; P=0.5559120396302443; Q=3586.890821629729; R=32282.01739466756; S=0.2403956607252189; T=0.2403956607252188 }
    6  = @{ E=3; G=34.88211266666666; H=104.646338; I=0.4324347083490296; J=0.4324347083490295; K=3; M=63.66262833333334; N=190.987885; O=0.3441717873742006; P=0.3441717873742006; Q=2220.686974179459; R=19986.18276761513; S=0.1488318264951267; T=0.1488318264951267 }
    7  = @{ E=3; G=34.88211266666666; H=104.646338; I=0.4324347083490296; J=0.4324347083490295; K=3; M=18.481835; N=55.445505; O=0.09991617299555507; P=0.09991617299555505; Q=644.6854507567433; R=5802.16905681069; S=0.04320722112868404; T=0.04320722112868403 }
    8  = @{ E=3; G=11.67052633333333; H=35.011579; I=0.1446799022599722; J=0.1446799022599721; K=3; M=102.8289443333334; N=308.486833; O=0.5559120396302444; P=0.5559120396302443; Q=1200.067902671034; R=10800.61112403931; S=0.08042929955884552; T=0.08042929955884549 }
    9  = @{ E=3; G=11.67052633333333; H=35.011579; I=0.1446799022599722; J=0.1446799022599721; K=3; M=63.66262833333334; N=190.987885; O=0.3441717873742006; P=0.3441717873742006; Q=742.9763804133795; R=6686.787423720415; S=0.04979474055793926; T=0.04979474055793925 }
    10 = @{ E=3; G=11.67052633333333; H=35.011579; I=0.1446799022599722; J=0.1446799022599721; K=3; M=18.481835; N=55.445505; O=0.09991617299555507; P=0.09991617299555505; Q=215.6927420558217; R=1941.234678502395; S=0.01445586214318738; T=0.01445586214318737 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
